# Automatic update of files.
# Applies:
#   1. "Förändrad" date (column C) bumped from 45184 to 45186 for all existing data rows (2-217)
#   2. Hyperlink formulas in rows 2-9 (columns S,T,U,V,W,X,Y) get a friendly-text second
#      argument equal to the case/"Beteckning" identifier (column A) of that row.
#   3. Row 217 gains an explicit custom row height (15) matching the new row below it.
#   4. A brand new data row (218) is appended for case "A 43521-2023".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastExistingRow = 217

# --- 1. Bump the "Förändrad" (changed) date for every existing data row -------------
$changedRange = $ws.Range("C2:C" + $lastExistingRow)
$changedRange.Value = 45186

# --- 2. Add friendly display text to HYPERLINK formulas on rows 2-9 -----------------
for ($r = 2; $r -le 9; $r++) {
    $caseId = $ws.Cells.Item($r, 1).Value2   # column A = "Beteckning"

    # Hyperlink columns are S(19), T(20), U(21), V(22), W(23), X(24), Y(25)
    for ($col = 19; $col -le 25; $col++) {
        $cell = $ws.Cells.Item($r, $col)
        $f = $cell.Formula

        if ($f -ne $null -and $f.Length -gt 0 -and $f.ToUpper().StartsWith("=HYPERLINK(") -and -not $f.Contains(",")) {
            $newFormula = $f.Substring(0, $f.Length - 1) + ', "' + $caseId + '")'
            $cell.Formula = $newFormula
        }
    }
}

# --- 3. Give row 217 the same explicit row height as the newly appended row --------
$ws.Rows.Item($lastExistingRow).RowHeight = 15

# --- 4. Append the new row (218) ----------------------------------------------------
$newRow = $lastExistingRow + 1

$ws.Cells.Item($newRow, 1).Value = "A 43521-2023"          # A - Beteckning

$ws.Cells.Item($newRow, 2).Value = 45184                    # B - Datum
$ws.Cells.Item($newRow, 2).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item($newRow, 3).Value = 45186                    # C - Förändrad
$ws.Cells.Item($newRow, 3).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item($newRow, 4).Value = "ÖSTERGÖTLANDS LÄN"      # D - Län
$ws.Cells.Item($newRow, 5).Value = "ÖDESHÖG"                # E - Kommun

$ws.Cells.Item($newRow, 7).Value = 1.1                      # G - Area (ha)

# H..Q (8..17) are all zero counts
for ($col = 8; $col -le 17; $col++) {
    $ws.Cells.Item($newRow, $col).Value = 0
}

# R - Artnamn: stays blank but keeps the wrap-text style used throughout column R
$ws.Cells.Item($newRow, 18).WrapText = $true
